$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, Coin (B), Link (C), Price (D), Volume(1h) (E)
$rows = @(
    @(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "69.113.59", "  -0.55%  "),
    @(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "2.469.22", "  -0.98%  "),
    @(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.00", "  +0.03%  "),
    @(5, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "559.86", "  -1.66%  "),
    @(6, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "163.11", "  -1.78%  "),
    @(7, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.00", "  +0.04%  "),
    @(8, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.504", "  -1.19%  "),
    @(9, "LidoStakedEther", "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth", "2.468.68", "  -0.99%  "),
    @(10, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.152", "  -4.71%  "),
    @(11, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.165", "  -0.60%  "),
    @(12, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.335", "  -3.59%  "),
    @(13, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "4.82", "  -1.13%  "),
    @(14, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "2.924.02", "  -0.93%  "),
    @(15, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "68.869.79", "  -0.75%  "),
    @(16, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.0000169", "  -2.90%  "),
    @(17, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "23.61", "  -2.29%  "),
    @(18, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "2.481.76", "  -0.30%  "),
    @(19, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "10.81", "  -3.47%  "),
    @(20, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "342.57", "  -2.95%  "),
    @(21, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "7.09", "  -4.84%  "),
    @(22, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "3.81", "  -2.54%  "),
    @(23, "SuiNetwork", "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui", "1.93", "  +1.18%  "),
    @(24, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.00", "  +0.05%  "),
    @(25, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "67.19", "  -2.95%  "),
    @(26, "NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "3.69", "  -2.68%  "),
    @(27, "WrappedeETH", "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth", "2.597.61", "  -0.92%  "),
    @(28, "Binance-PegBSC-USD", "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd", "1.01", "  +0.50%  "),
    @(29, "Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "8.21", "  -5.17%  "),
    @(30, "PEPE", "https://coinranking.com/coin/03WI8NQPF+pepe-pepe", "0.0₃0818", "  -6.14%  "),
    @(31, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "7.19", "  -4.64%  "),
    @(32, "Bittensor", "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao", "440.85", "  +0.55%  "),
    @(33, "FirstDigitalUSD", "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd", "1.00", "  +0.04%  "),
    @(34, "Fetch.AI", "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet", "1.14", "  -4.12%  "),
    @(35, "PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "1.62", "  -5.27%  "),
    @(36, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "156.16", "  +2.07%  "),
    @(37, "WhiteBITCoin", "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt", "19.05", "  -0.05%  "),
    @(38, "USDe", "https://coinranking.com/coin/exbfr2U-0+usde-usde", "1.00", "  -0.08%  "),
    @(39, "Kaspa", "https://coinranking.com/coin/V8GxkwWow+kaspa-kas", "0.109", "  -4.13%  "),
    @(40, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "17.92", "  -1.20%  "),
    @(41, "PolygonEcosystemToken", "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol", "0.302", "  -3.65%  "),
    @(42, "RenderToken", "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render", "4.46", "  -2.80%  "),
    @(43, "OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "37.44", "  -0.97%  "),
    @(44, "Stacks", "https://coinranking.com/coin/mMPrMcB7+stacks-stx", "1.48", "  -5.88%  "),
    @(45, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "1.10", "  +2.76%  "),
    @(46, "dogwifhat", "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif", "2.08", "  -4.16%  "),
    @(47, "Aave", "https://coinranking.com/coin/ixgUfzmLR+aave-aave", "133.35", "  -4.17%  "),
    @(48, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "3.36", "  -2.13%  "),
    @(49, "Cronos", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro", "0.0719", "  -0.49%  "),
    @(50, "ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "0.483", "  -4.23%  "),
    @(51, "Mantle", "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt", "0.563", "  -1.71%  ")
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 2).Value = $r[1]
    $ws.Cells.Item($rowNum, 3).Value = $r[2]

    # Price column: force text storage so numeric-looking strings
    # (e.g. "1.00", "163.11") are not auto-converted to numbers,
    # then restore the default "Normal" style so no extra number
    # format is left attached to the cell.
    $priceCell = $ws.Cells.Item($rowNum, 4)
    $priceCell.NumberFormat = "@"
    $priceCell.Value = $r[3]
    $priceCell.Style = "Normal"

    $ws.Cells.Item($rowNum, 5).Value = $r[4]
}

Write-Host "Done updating cryptos list"